# DOMA-4473: add an "isVerified" column to the contacts export template.
#
# The template stores its column headers/placeholders as plain text in
# row 1 (i18n labels) and rows 2-3 (the {d.contacts[...]} field tokens for
# the "current" and "next" row of the generated report). A new
# "isVerified" field is introduced as the FIRST field of each group, which
# pushes every existing field one column to the right (name -> B,
# address -> C, ... role -> H) and a brand new column H is populated with
# what used to be in G (role).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new column H the same look (font/fill/border) as the rest of
# its row by copying the formatting from the last existing column (G)
# before any values are written into it.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("G3").Copy()
$ws.Range("H3").PasteSpecial(-4122)

# Row 1 - {d.i18n.*} header labels. Shift each label one column to the
# right (H <- G <- F <- ... <- B <- A), then seed A with the new field.
$ws.Range("H1").Value = "{d.i18n.role}"
$ws.Range("G1").Value = "{d.i18n.email}"
$ws.Range("F1").Value = "{d.i18n.phone}"
$ws.Range("E1").Value = "{d.i18n.unitType}"
$ws.Range("D1").Value = "{d.i18n.unitName}"
$ws.Range("C1").Value = "{d.i18n.address}"
$ws.Range("B1").Value = "{d.i18n.name}"
$ws.Range("A1").Value = "{d.i18n.isVerified}"

# Row 2 - {d.contacts[i].*} tokens for the current contact row.
$ws.Range("H2").Value = "{d.contacts[i].role}"
$ws.Range("G2").Value = "{d.contacts[i].email}"
$ws.Range("F2").Value = "{d.contacts[i].phone}"
$ws.Range("E2").Value = "{d.contacts[I].unitType}"
$ws.Range("D2").Value = "{d.contacts[i].unitName}"
$ws.Range("C2").Value = "{d.contacts[i].address}"
$ws.Range("B2").Value = "{d.contacts[i].name}"
$ws.Range("A2").Value = "{d.contacts[i].isVerified}"

# Row 3 - {d.contacts[i+1].*} tokens for the next contact row.
$ws.Range("H3").Value = "{d.contacts[i+1].role}"
$ws.Range("G3").Value = "{d.contacts[i+1].email}"
$ws.Range("F3").Value = "{d.contacts[i+1].phone}"
$ws.Range("E3").Value = "{d.contacts[I+1].unitType}"
$ws.Range("D3").Value = "{d.contacts[i+1].unitName}"
$ws.Range("C3").Value = "{d.contacts[i+1].address}"
$ws.Range("B3").Value = "{d.contacts[i+1].name}"
$ws.Range("A3").Value = "{d.contacts[i+1].isVerified}"

Write-Host "Added isVerified column to contacts export template"
